$d = $word.ActiveDocument

# The "Micro results" row's value cell (table 1, column 2) currently holds a
# single paragraph whose lines are joined with <w:br/>. Replace it with one
# paragraph per result line (plus a leading blank paragraph), each run in
# blue Times New Roman 10pt, matching the refreshed micro-results summary.
$table = $d.Tables.Item(1)
$targetCell = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $labelText = $table.Cell($i, 1).Range.Text
    if ($labelText -like "*Micro results*") {
        $targetCell = $table.Cell($i, 2)
        break
    }
}

if ($targetCell -eq $null) {
    throw "Could not find 'Micro results' row"
}

# NOTE: the fragment passed to InsertXML must be the bare paragraph-level
# WordprocessingML content (no enclosing <w:tc>/<w:body>), otherwise Word
# replaces the cell's contents with nothing instead of the new paragraphs.
$resultsXml = @'
<w:p><w:r/></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>21/04 – SARS CORONAVIRUS–2 PCR – Negative</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>21/04 – SARS–CoV–2 RNA – Negative</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>19/04 – BLC – PERIPHERAL–RIGHT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>19/04 – UC – **No clear Result**</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic recommendations without pathogen details.</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>19/04 – UC – CATHETER SPECIMEN URINE (CSU) NO SIGNIFICANT GROWTH</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>--------Previous result (1 year)--------</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>17/01 – VZV IgG – Positive</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>17/01 – EBNA IgG – Positive</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>17/01 – CMV IgG – Positive</w:t></w:r></w:p>
'@

$targetCell.Range.InsertXML($resultsXml)
